# Applies the NATMI recompute for Hgf -> St14 ligand-receptor pairs
# (ECs/FAPs/M2/sCs source x ECs/FAPs/M2/sCs target = 16 rows, rows 2-17).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ECs", "Hgf", "St14", "ECs", 2, 0.6666666666666666, 0.1088396666666667, 0.326519, 0.002750770615347974, 0.002750770615347974, 3, 1, 0.7971666666666667, 2.3915, 0.3308932235309289, 0.3308932235309289, 0.08676335427777779, 0.7808701885000001, 0.0009102113561066478, 0.0009102113561066479),
    @("ECs", "Hgf", "St14", "FAPs", 2, 0.6666666666666666, 0.1088396666666667, 0.326519, 0.002750770615347974, 0.002750770615347974, 3, 1, 0.679891, 2.039673, 0.282213662521012, 0.282213662521012, 0.07399910980966667, 0.665991988287, 0.0007763050501125296, 0.0007763050501125297),
    @("ECs", "Hgf", "St14", "M2", 2, 0.6666666666666666, 0.1088396666666667, 0.326519, 0.002750770615347974, 0.002750770615347974, 3, 1, 0.76147, 2.28441, 0.3160760145374406, 0.3160760145374406, 0.08287814097666667, 0.7459032687899999, 0.0008694526130058906, 0.0008694526130058906),
    @("ECs", "Hgf", "St14", "sCs", 2, 0.6666666666666666, 0.1088396666666667, 0.326519, 0.002750770615347974, 0.002750770615347974, 3, 1, 0.170608, 0.5118240000000001, 0.07081709941061849, 0.07081709941061851, 0.01856891785066667, 0.167120260656, 0.0001948015961229057, 0.0001948015961229057),
    @("FAPs", "Hgf", "St14", "ECs", 3, 1, 15.07419333333333, 45.22258, 0.3809791902285103, 0.3809791902285103, 3, 1, 0.7971666666666667, 2.3915, 0.3308932235309289, 0.3308932235309289, 12.01664445222222, 108.14980007, 0.1260634323529148, 0.1260634323529148),
    @("FAPs", "Hgf", "St14", "FAPs", 3, 1, 15.07419333333333, 45.22258, 0.3809791902285103, 0.3809791902285103, 3, 1, 0.679891, 2.039673, 0.282213662521012, 0.282213662521012, 10.24880837959333, 92.23927541634001, 0.1075175326186773, 0.1075175326186773),
    @("FAPs", "Hgf", "St14", "M2", 3, 1, 15.07419333333333, 45.22258, 0.3809791902285103, 0.3809791902285103, 3, 1, 0.76147, 2.28441, 0.3160760145374406, 0.3160760145374406, 11.47854599753333, 103.3069139778, 0.120418384069129, 0.120418384069129),
    @("FAPs", "Hgf", "St14", "sCs", 3, 1, 15.07419333333333, 45.22258, 0.3809791902285103, 0.3809791902285103, 3, 1, 0.170608, 0.5118240000000001, 0.07081709941061849, 0.07081709941061851, 2.571777976213334, 23.14600178592, 0.02697984118778935, 0.02697984118778935),
    @("M2", "Hgf", "St14", "ECs", 3, 1, 23.69325166666667, 71.07975500000001, 0.5988138558556568, 0.5988138558556569, 3, 1, 0.7971666666666667, 2.3915, 0.3308932235309289, 0.3308932235309289, 18.88747045361111, 169.9872340825, 0.1981434470590633, 0.1981434470590633),
    @("M2", "Hgf", "St14", "FAPs", 3, 1, 23.69325166666667, 71.07975500000001, 0.5988138558556568, 0.5988138558556569, 3, 1, 0.679891, 2.039673, 0.282213662521012, 0.282213662521012, 16.10882856890167, 144.979457120115, 0.1689934514293543, 0.1689934514293543),
    @("M2", "Hgf", "St14", "M2", 3, 1, 23.69325166666667, 71.07975500000001, 0.5988138558556568, 0.5988138558556569, 3, 1, 0.76147, 2.28441, 0.3160760145374406, 0.3160760145374406, 18.04170034661667, 162.37530311955, 0.1892706970086535, 0.1892706970086535),
    @("M2", "Hgf", "St14", "sCs", 3, 1, 23.69325166666667, 71.07975500000001, 0.5988138558556568, 0.5988138558556569, 3, 1, 0.170608, 0.5118240000000001, 0.07081709941061849, 0.07081709941061851, 4.042258280346667, 36.38032452312001, 0.04240626035858582, 0.04240626035858584),
    @("sCs", "Hgf", "St14", "ECs", 3, 1, 0.6906883333333335, 2.072065, 0.01745618330048481, 0.01745618330048481, 3, 1, 0.7971666666666667, 2.3915, 0.3308932235309289, 0.3308932235309289, 0.550593716388889, 4.955343447500001, 0.005776132762844188, 0.005776132762844188),
    @("sCs", "Hgf", "St14", "FAPs", 3, 1, 0.6906883333333335, 2.072065, 0.01745618330048481, 0.01745618330048481, 3, 1, 0.679891, 2.039673, 0.282213662521012, 0.282213662521012, 0.4695927816383335, 4.226335034745, 0.004926373422867947, 0.004926373422867947),
    @("sCs", "Hgf", "St14", "M2", 3, 1, 0.6906883333333335, 2.072065, 0.01745618330048481, 0.01745618330048481, 3, 1, 0.76147, 2.28441, 0.3160760145374406, 0.3160760145374406, 0.5259384451833334, 4.73344600665, 0.005517480846652265, 0.005517480846652264),
    @("sCs", "Hgf", "St14", "sCs", 3, 1, 0.6906883333333335, 2.072065, 0.01745618330048481, 0.01745618330048481, 3, 1, 0.170608, 0.5118240000000001, 0.07081709941061849, 0.07081709941061851, 0.1178369551733334, 1.06053259656, 0.001236196268120411, 0.001236196268120411)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
